# Auto-generated edit script for Maduin_Profits workbook update
# Updates currentAveragePrice / LevePrice / LeveProfit columns across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5340.846
$ws.Range("I64").Value = 4571.4287
$ws.Range("J64").Value = 6238.5
$ws.Range("K64").Value = 4571.4287
$ws.Range("L64").Value = 6238.5
$ws.Range("M64").Value = -4323.4287
$ws.Range("N64").Value = -6734.5
$ws.Range("H67").Value = 5340.846
$ws.Range("I67").Value = 4571.4287
$ws.Range("J67").Value = 6238.5
$ws.Range("K67").Value = 4571.4287
$ws.Range("L67").Value = 6238.5
$ws.Range("M67").Value = -3713.4287
$ws.Range("N67").Value = -7954.5
$ws.Range("H100").Value = 2368
$ws.Range("I100").Value = 2368
$ws.Range("K100").Value = 2368
$ws.Range("M100").Value = -1827

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3097.2144
$ws.Range("I2").Value = 1057
$ws.Range("K2").Value = 1057
$ws.Range("M2").Value = -944
$ws.Range("H4").Value = 766.3333
$ws.Range("I4").Value = 649.5
$ws.Range("K4").Value = 649.5
$ws.Range("M4").Value = -533.5
$ws.Range("H14").Value = 1606
$ws.Range("I14").Value = 1606
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1606
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1431
$ws.Range("N14").Value = $null
$ws.Range("H32").Value = 2463
$ws.Range("I32").Value = 2262.32
$ws.Range("J32").Value = 4971.5
$ws.Range("K32").Value = 2262.32
$ws.Range("L32").Value = 4971.5
$ws.Range("M32").Value = -1975.32
$ws.Range("N32").Value = -5545.5
$ws.Range("H61").Value = 2500
$ws.Range("I61").Value = 2500
$ws.Range("K61").Value = 2500
$ws.Range("M61").Value = -2288
$ws.Range("H102").Value = 2642.7144
$ws.Range("I102").Value = 2833.1667
$ws.Range("K102").Value = 2833.1667
$ws.Range("M102").Value = -1211.1667
$ws.Range("H116").Value = 3097.2144
$ws.Range("I116").Value = 1057
$ws.Range("K116").Value = 1057
$ws.Range("M116").Value = 1237
$ws.Range("H136").Value = 2500
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -4950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3097.2144
$ws.Range("I3").Value = 1057
$ws.Range("K3").Value = 1057
$ws.Range("M3").Value = -943
$ws.Range("H20").Value = 1175.1
$ws.Range("I20").Value = 536
$ws.Range("J20").Value = 2133.75
$ws.Range("K20").Value = 536
$ws.Range("L20").Value = 2133.75
$ws.Range("M20").Value = -289
$ws.Range("N20").Value = -2627.75
$ws.Range("H22").Value = 514.25
$ws.Range("I22").Value = 608.7778
$ws.Range("J22").Value = 230.66667
$ws.Range("K22").Value = 608.7778
$ws.Range("L22").Value = 230.66667
$ws.Range("M22").Value = -435.7778
$ws.Range("N22").Value = -576.6666700000001
$ws.Range("H64").Value = 1017.4
$ws.Range("I64").Value = 999
$ws.Range("J64").Value = 1022
$ws.Range("K64").Value = 999
$ws.Range("L64").Value = 1022
$ws.Range("M64").Value = -774
$ws.Range("N64").Value = -1472
$ws.Range("H67").Value = 1017.4
$ws.Range("I67").Value = 999
$ws.Range("J67").Value = 1022
$ws.Range("K67").Value = 999
$ws.Range("L67").Value = 1022
$ws.Range("M67").Value = -219
$ws.Range("N67").Value = -2582
$ws.Range("H107").Value = 1858.3334
$ws.Range("J107").Value = 1900
$ws.Range("L107").Value = 1900
$ws.Range("N107").Value = -5740

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2862
$ws.Range("J31").Value = 3150
$ws.Range("L31").Value = 3150
$ws.Range("N31").Value = -3740
$ws.Range("H34").Value = 2862
$ws.Range("J34").Value = 3150
$ws.Range("L34").Value = 3150
$ws.Range("N34").Value = -3554

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 2000
$ws.Range("I76").Value = 2000
$ws.Range("K76").Value = 6000
$ws.Range("M76").Value = -5617
$ws.Range("H79").Value = 2000
$ws.Range("I79").Value = 2000
$ws.Range("K79").Value = 6000
$ws.Range("M79").Value = -4674
$ws.Range("H107").Value = 768
$ws.Range("I107").Value = 503
$ws.Range("K107").Value = 1509
$ws.Range("M107").Value = 411
$ws.Range("H131").Value = 977.3
$ws.Range("I131").Value = 799
$ws.Range("J131").Value = 1021.875
$ws.Range("K131").Value = 2397
$ws.Range("L131").Value = 3065.625
$ws.Range("M131").Value = 2643
$ws.Range("N131").Value = -13145.625
$ws.Range("H132").Value = 2320
$ws.Range("I132").Value = 2189.875
$ws.Range("J132").Value = 2667
$ws.Range("K132").Value = 19708.875
$ws.Range("L132").Value = 24003
$ws.Range("M132").Value = -17178.875
$ws.Range("N132").Value = -29063
$ws.Range("H134").Value = 2874.75
$ws.Range("I134").Value = 2874.75
$ws.Range("K134").Value = 8624.25
$ws.Range("M134").Value = -3554.25
$ws.Range("H139").Value = 2539.5
$ws.Range("I139").Value = 3398.75
$ws.Range("J139").Value = 1966.6666
$ws.Range("K139").Value = 10196.25
$ws.Range("L139").Value = 5899.9998
$ws.Range("M139").Value = -5056.25
$ws.Range("N139").Value = -16179.9998
$ws.Range("H140").Value = 1102.7142
$ws.Range("I140").Value = 953.1667
$ws.Range("K140").Value = 2859.5001
$ws.Range("M140").Value = 2320.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null
$ws.Range("H15").Value = 24999.666
$ws.Range("J15").Value = 24999.666
$ws.Range("L15").Value = 24999.666
$ws.Range("N15").Value = -25575.666
$ws.Range("H81").Value = 24999.666
$ws.Range("J81").Value = 24999.666
$ws.Range("L81").Value = 24999.666
$ws.Range("N81").Value = -26995.666
$ws.Range("H84").Value = 24999.666
$ws.Range("J84").Value = 24999.666
$ws.Range("L84").Value = 74998.998
$ws.Range("N84").Value = -84982.998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 591.5714
$ws.Range("J55").Value = 591.5714
$ws.Range("L55").Value = 591.5714
$ws.Range("N55").Value = -937.5714
$ws.Range("H68").Value = 2333
$ws.Range("J68").Value = 999
$ws.Range("L68").Value = 999
$ws.Range("N68").Value = -2497
$ws.Range("H71").Value = 2333
$ws.Range("J71").Value = 999
$ws.Range("L71").Value = 4995
$ws.Range("N71").Value = -12483
$ws.Range("H132").Value = 4451.9287
$ws.Range("I132").Value = 4025.1538
$ws.Range("K132").Value = 12075.4614
$ws.Range("M132").Value = -9545.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1006.8
$ws.Range("I107").Value = 1021.125
$ws.Range("J107").Value = 949.5
$ws.Range("K107").Value = 3063.375
$ws.Range("L107").Value = 2848.5
$ws.Range("M107").Value = -1143.375
$ws.Range("N107").Value = -6688.5
$ws.Range("H132").Value = 1283.25
$ws.Range("I132").Value = 1283.25
$ws.Range("K132").Value = 3849.75
$ws.Range("M132").Value = -1319.75
